# Commiting the acc to legal agreement flow
# Adds a new "AcctoLegalAgreement" test case row to the TestCases sheet,
# switches the previously-active "LeadtoQuote" row to RunMode "No",
# and sets the new row's RunMode to "Yes" so only the new flow runs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing last data row (row 9) was "LeadtoQuote" with RunMode "Yes".
# It now becomes "No" since the newly added test case takes over as the
# active one to run.
$ws.Range("C9").Value = "No"

# New row 10: the added test case and its RunMode ("Yes").
$ws.Range("A10").Value = "AcctoLegalAgreement"
$ws.Range("C10").Value = "Yes"

# Duplicate the formatting of C10 (column C body style) down onto the two
# new blank trailer rows (11 and 12), matching the sheet's existing
# pattern of a formatted-but-empty column C cell trailing the data.
$ws.Range("C10").Copy()
$ws.Range("C11:C12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match row heights used throughout the sheet.
$ws.Rows.Item(11).RowHeight = 15.75
$ws.Rows.Item(12).RowHeight = 15.75

# Move the active selection to A12, as recorded in the saved view state.
$ws.Range("A12").Select()
